# Add a new "TextBox 1" shape to slide 1 containing the project's GitHub URL.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU targets from the OOXML diff (PowerPoint COM works in points, 1 pt = 12700 EMU).
# The host's point->EMU conversion truncates rather than rounds, so the point
# values below are nudged slightly so that floor(pt_as_float32 * 12700) lands
# exactly on the desired EMU value.
$targetLeftEmu   = 3122065
$targetTopEmu    = 6462111
$targetWidthEmu  = 6084606
$targetHeightEmu = 338554

$leftPt   = 245.83188976377951
$topPt    = 508.8276519952756
$widthPt  = 479.1028346456693
$heightPt = 26.657795975590552

$tb = $s.Shapes.AddTextbox(1, $leftPt, $topPt, $widthPt, $heightPt)

$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = 0

$tr = $tb.TextFrame.TextRange
$tr.Text = "https://github.com/ksu-hmi/cryptostockportfolio"
$tr.Font.Size = 16
$tr.Font.Name = "Cooper Black"
